# Update view-count figures (column F) on the "展览" and "全部类型" sheets
# to reflect the refreshed data snapshot (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15662
$ws1.Range("F8").Value  = 706
$ws1.Range("F9").Value  = 15427
$ws1.Range("F11").Value = 9024
$ws1.Range("F12").Value = 383
$ws1.Range("F18").Value = 199
$ws1.Range("F21").Value = 550
$ws1.Range("F27").Value = 16
$ws1.Range("F33").Value = 63
$ws1.Range("F36").Value = 324
$ws1.Range("F37").Value = 453
$ws1.Range("F39").Value = 5552

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15662
$ws4.Range("F8").Value  = 706
$ws4.Range("F9").Value  = 15427
$ws4.Range("F11").Value = 9024
$ws4.Range("F12").Value = 383
$ws4.Range("F18").Value = 199
$ws4.Range("F21").Value = 550
$ws4.Range("F27").Value = 16
$ws4.Range("F35").Value = 63
$ws4.Range("F38").Value = 324
$ws4.Range("F39").Value = 453
$ws4.Range("F41").Value = 5552
